# Apply the Brachiopoda Materials-sheet review-comment fixes:
#  - insert suborder / infraorder / superfamily columns between order and family
#  - point scientificName at the new summary.taxonName ARPHA placeholder
#  - populate scientificNameAuthorship from summary.Author
#  - populate the (previously empty) eventTime column with a time placeholder

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Insert three new blank columns right after "order" (column AQ), i.e. where
# "family" currently starts (AR:AT), pushing family/genus/subgenus/... right.
$ws.Range("AR1:AT1").EntireColumn.Insert()

# New header row (row 1) labels for the inserted columns.
$ws.Range("AR1").Value = "suborder"
$ws.Range("AS1").Value = "infraorder"
$ws.Range("AT1").Value = "superfamily"

# New data row (row 2) ARPHA template placeholders for the inserted columns.
$ws.Range("AR2").Value = "`${suborder}"
$ws.Range("AS2").Value = "`${infraorder}"
$ws.Range("AT2").Value = "`${superfamily}"

# scientificName (column AG) now resolves from the iNaturalist summary object.
$ws.Range("AG2").Value = "`${summary.taxonName}"

# scientificNameAuthorship (was column AY, now shifted to BB after the insert)
# picks up the authorship straight from the summary object too.
$ws.Range("BB2").Value = "`${summary.Author}"

# eventTime (was column DX, now shifted to EA after the insert) gets a time
# template alongside the existing eventDate (YYYY-MM-DD) template.
$ws.Range("EA2").Value = "!Date:HH:mm:ss"
